# Auto-upload VRF Excel file
# Adds a new "cat" worksheet at the end of the workbook with a single
# header row describing outdoor/indoor unit model, quantity and serial
# number columns.

$wb = $excel.ActiveWorkbook

# Remember which sheet was active so we can restore the selection once
# the new sheet has been appended (Worksheets.Add() activates the sheet
# it creates).
$originalActiveSheetName = $wb.ActiveSheet.Name

# Insert the new sheet after the current last sheet so it lands at the
# very end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "cat"

$headers = @("Outdoor Model", "Outdoor Quantity", "Outdoor Serial(s)", "Indoor Model", "Indoor Quantity", "Indoor Serial(s)")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$null = $ws.Range("A1").Select()

# Restore the original active sheet / selection.
$null = $wb.Worksheets.Item($originalActiveSheetName).Activate()
